$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New instrument / test-case identifiers (TC01/DNBSEQ-G400 -> TC11/Illumina NovaSeq 6000)
$neo4jFile = "TC11_CDS_Filter_InstrumentModel-Illumina NovaSeq 6000_Neo4jData.xlsx"
$webFile   = "TC11_CDS_Filter_InstrumentModel-Illumina NovaSeq 6000_WebData.xlsx"

# Query texts (instrument model updated to "Illumina NovaSeq 6000")
$participantsQuery = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina NovaSeq 6000']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p, s, collect(distinct samp.sample_id) as samp`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY ``Participant ID``LIMIT 100"

$samplesQuery = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina NovaSeq 6000']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN  `n coalesce(samp.sample_id, '') as ``Sample ID``,`n coalesce(p.participant_id,'') as ``Participant ID``,`n coalesce(s.study_name, '') as ``Study Name``,`n coalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(samp.sample_tumor_status,'') as ``Tumor``,`ncoalesce(samp.sample_type,'') as ``Analyte Type```nORDER By samp.sample_id LIMIT 100"

$filesListQuery = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina NovaSeq 6000']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN `n    coalesce(f.file_name, '') as ``File Name``,`n    coalesce(s.study_name, '') as ``Study Name``,`n    coalesce(s.phs_accession,'') as ``Accession``,`n    coalesce(p.participant_id,'') as ``Participant ID``,`n    coalesce(samp.sample_id, '') as ``Sample ID``,`n    coalesce(f.file_type, '') as ``File Type```nORDER By f.file_name LIMIT 100"

$filesStatQuery = "MATCH (f:file)`nMatch (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina NovaSeq 6000']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,f, s, collect(distinct samp.sample_id) as samp`nRETURN`ncount(distinct s) AS Studies,`ncount(distinct p) AS Participants,`ncount(distinct samp) AS Samples,`ncount(distinct f) AS Files"

$filesStatQueryVariant = "MATCH (f:file)`nMatch (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina NovaSeq 6000']MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,f, s, collect(distinct samp.sample_id) as samp`nRETURN`ncount(distinct s) AS Studies,`ncount(distinct p) AS Participants,`ncount(distinct samp) AS Samples,`ncount(distinct f) AS Files"

# Set cells in the same order the original authors did (file names first, then
# each tab's query, then the StatQuery column) so that newly-created shared
# strings land in the same order as the target workbook.
$ws.Range("D2").Value = $neo4jFile
$ws.Range("D3").Value = $neo4jFile
$ws.Range("D4").Value = $neo4jFile
$ws.Range("E2").Value = $webFile
$ws.Range("E3").Value = $webFile
$ws.Range("E4").Value = $webFile

# Row 2 - ParticipantsTab
$ws.Range("B2").Value = $participantsQuery
# Row 3 - SamplesTab
$ws.Range("B3").Value = $samplesQuery
# Row 4 - FilesTab
$ws.Range("B4").Value = $filesListQuery

# StatQuery column (well-formed query reused by rows 2 & 3)
$ws.Range("C2").Value = $filesStatQuery
$ws.Range("C3").Value = $filesStatQuery
# Row 4's StatQuery is a slightly malformed variant (missing a newline)
$ws.Range("C4").Value = $filesStatQueryVariant

# Column D widened to fit the new (longer) content
$ws.Columns.Item(4).ColumnWidth = 98

# Selection moved to C4
$ws.Range("C4").Select()
